$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet Hoja1 (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.94 = 36800.71 pesos`n✅ 36800.71 pesos = 8.88 = 958.63 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 111.9
$wsTasas.Range("O10").Value = 4118
$wsTasas.Range("N12").Value = 4146
$wsTasas.Range("O12").Value = 108
